$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the default (unformatted) cell style so we can restore it
# after temporarily marking risky cells as Text, to avoid corrupting
# numeric-looking strings (e.g. "1.002") into actual numbers.
$defaultStyle = $ws.Range("A1").Style

function Set-TextCellValue($rangeRef, $val) {
    $cell = $ws.Range($rangeRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $defaultStyle
}

$ws.Range("D2").Value = '27.557.61'
$ws.Range("E2").Value = '  +5.77%  '
$ws.Range("D3").Value = '1.814.25'
$ws.Range("E3").Value = '  +5.87%  '
Set-TextCellValue "D4" '1.002'
$ws.Range("E4").Value = '  +0.66%  '
Set-TextCellValue "D5" '343.92'
$ws.Range("E5").Value = '  +3.66%  '
Set-TextCellValue "D6" '0.9995'
$ws.Range("E6").Value = '  +0.15%  '
Set-TextCellValue "D7" '0.3838'
$ws.Range("E7").Value = '  +4.15%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCellValue "D8" '0.3516'
$ws.Range("E8").Value = '  +5.80%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCellValue "D9" '49.87'
$ws.Range("E9").Value = '  +2.96%  '
$ws.Range("E10").Value = '  +4.66%  '
Set-TextCellValue "D11" '0.07772'
$ws.Range("E11").Value = '  +3.75%  '
Set-TextCellValue "D12" '1.001'
$ws.Range("E12").Value = '  +0.80%  '
Set-TextCellValue "D13" '22.43'
$ws.Range("E13").Value = '  +11.56%  '
Set-TextCellValue "D14" '6.618'
$ws.Range("E14").Value = '  +6.34%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCellValue "D15" '7.214'
$ws.Range("E15").Value = '  +4.66%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.813.95'
$ws.Range("E16").Value = '  +6.50%  '
Set-TextCellValue "D17" '0.00001129'
$ws.Range("E17").Value = '  +5.08%  '
Set-TextCellValue "D18" '0.06717'
$ws.Range("E18").Value = '  +1.20%  '
Set-TextCellValue "D19" '86.58'
$ws.Range("E20").Value = '  +0.35%  '
Set-TextCellValue "D21" '17.77'
$ws.Range("E21").Value = '  +9.31%  '
Set-TextCellValue "D22" '6.539'
$ws.Range("E22").Value = '  +7.85%  '
Set-TextCellValue "D23" '13.15'
$ws.Range("E23").Value = '  +0.86%  '
$ws.Range("D24").Value = '27.546.36'
$ws.Range("E24").Value = '  +6.31%  '
Set-TextCellValue "D25" '2.468'
$ws.Range("E25").Value = '  +0.21%  '
Set-TextCellValue "D26" '2.682'
$ws.Range("E26").Value = '  +7.69%  '
Set-TextCellValue "D27" '22.14'
$ws.Range("E27").Value = '  +15.19%  '
Set-TextCellValue "D28" '1.498'
$ws.Range("E28").Value = '  +14.47%  '
Set-TextCellValue "D29" '153.94'
$ws.Range("E29").Value = '  +2.89%  '
$ws.Range("D30").Value = '2.017.27'
$ws.Range("E30").Value = '  +6.72%  '
Set-TextCellValue "D31" '136.64'
$ws.Range("E31").Value = '  +6.51%  '
Set-TextCellValue "D32" '6.388'
$ws.Range("E32").Value = '  +7.14%  '
Set-TextCellValue "D33" '4.086'
$ws.Range("E33").Value = '  -0.52%  '
Set-TextCellValue "D34" '13.98'
$ws.Range("E34").Value = '  +8.30%  '
Set-TextCellValue "D35" '0.08830'
$ws.Range("E35").Value = '  +3.79%  '
Set-TextCellValue "D36" '1.722'
$ws.Range("E36").Value = '  +0.03%  '
Set-TextCellValue "D37" '5.647'
$ws.Range("E37").Value = '  +5.54%  '
Set-TextCellValue "D38" '0.7118'
Set-TextCellValue "D39" '0.06551'
$ws.Range("E39").Value = '  +5.51%  '
Set-TextCellValue "D40" '0.2270'
$ws.Range("E40").Value = '  +7.04%  '
Set-TextCellValue "D41" '0.02418'
$ws.Range("E41").Value = '  +6.17%  '
$ws.Range("E42").Value = '  +5.60%  '
Set-TextCellValue "D43" '1.286'
$ws.Range("E43").Value = '  +0.81%  '
Set-TextCellValue "D44" '14.88'
$ws.Range("E44").Value = '  +2.10%  '
Set-TextCellValue "D45" '0.6650'
$ws.Range("E45").Value = '  +13.50%  '
Set-TextCellValue "D46" '0.9996'
$ws.Range("E46").Value = '  +0.14%  '
Set-TextCellValue "D47" '4.036'
$ws.Range("E47").Value = '  +5.18%  '
Set-TextCellValue "D48" '2.184'
$ws.Range("E48").Value = '  +8.80%  '
Set-TextCellValue "D49" '133.07'
$ws.Range("E49").Value = '  +4.69%  '
$ws.Range("E50").Value = '  +1.76%  '
Set-TextCellValue "D51" '80.79'
$ws.Range("E51").Value = '  +5.31%  '
